{"js": "// 1. Every table gets a \"start\" (logical-left) table justification:\n//    <w:tblPr> ... <w:jc w:val=\"start\"/> ... </w:tblPr>\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\nfor (const table of tables.items) {\n  table.alignment = \"start\";\n}\nawait context.sync();\n\n// 2. New custom paragraph style \"Abstract Title\", based on Normal,\n//    followed by the existing \"Abstract\" style.\ncontext.document.addStyle(\"Abstract Title\", Word.StyleType.paragraph);\nawait context.sync();\n\n// addStyle's return value isn't reliably bound in every host, so re-fetch\n// the style by name before touching its properties.\nconst abstractTitle = context.document.getStyles().getByNameOrNullObject(\"Abstract Title\");\nawait context.sync();\n\nabstractTitle.baseStyle = \"Normal\";\nabstractTitle.nextParagraphStyle = \"Abstract\";\nabstractTitle.quickStyle = true;\nabstractTitle.paragraphFormat.keepWithNext = true;\nabstractTitle.paragraphFormat.keepTogether = true;\nabstractTitle.paragraphFormat.alignment = Word.Alignment.centered;\nabstractTitle.paragraphFormat.spaceAfter = 0;\nabstractTitle.paragraphFormat.spaceBefore = 15; // 300 twips = 15pt\nabstractTitle.font.size = 10;\nabstractTitle.font.sizeBidirectional = 10;\nabstractTitle.font.bold = true;\nabstractTitle.font.color = \"#345A8A\";\nawait context.sync();\n\n// 3. \"Abstract\" style: space-before shrinks from 300 to 100 twips (15pt -> 5pt).\nconst abstractStyle = context.document.getStyles().getByNameOrNullObject(\"Abstract\");\nawait context.sync();\nabstractStyle.paragraphFormat.spaceBefore = 5;\nawait context.sync();\n\n// 4. \"ImportTok\" character style: green + bold.\nconst importTok = context.document.getStyles().getByNameOrNullObject(\"ImportTok\");\nawait context.sync();\nimportTok.font.color = \"#008000\";\nimportTok.font.bold = true;\nawait context.sync();\n\n// 5. \"BuiltInTok\" character style: green.\nconst builtInTok = context.document.getStyles().getByNameOrNullObject(\"BuiltInTok\");\nawait context.sync();\nbuiltInTok.font.color = \"#008000\";\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1. Every table gets a \"start\" (logical-left) table justification:\n#    <w:tblPr> ... <w:jc w:val=\"start\"/> ... </w:tblPr>\nforeach ($t in $d.Tables) {\n    $t.Alignment = \"start\"\n}\n\n# 2. New custom paragraph style \"Abstract Title\", based on Normal,\n#    followed by the existing \"Abstract\" style.\n$abstractTitle = $d.Styles.Add(\"Abstract Title\", 1)\n$abstractTitle.BaseStyle = \"Normal\"\n$abstractTitle.NextParagraphStyle = \"Abstract\"\n$abstractTitle.QuickStyle = $true\n$abstractTitle.ParagraphFormat.KeepWithNext = $true\n$abstractTitle.ParagraphFormat.KeepTogether = $true\n$abstractTitle.ParagraphFormat.Alignment = \"center\"\n$abstractTitle.ParagraphFormat.SpaceAfter = 0\n$abstractTitle.ParagraphFormat.SpaceBefore = 15\n$abstractTitle.Font.Size = 10\n$abstractTitle.Font.SizeBi = 10\n$abstractTitle.Font.Bold = $true\n$abstractTitle.Font.Color = \"#345A8A\"\n\n# 3. \"Abstract\" style: space-before shrinks from 300 to 100 twips (15pt -> 5pt).\n$abstract = $d.Styles.Item(\"Abstract\")\n$abstract.ParagraphFormat.SpaceBefore = 5\n\n# 4. \"ImportTok\" character style: green + bold.\n$importTok = $d.Styles.Item(\"ImportTok\")\n$importTok.Font.Color = \"#008000\"\n$importTok.Font.Bold = $true\n\n# 5. \"BuiltInTok\" character style: green.\n$builtInTok = $d.Styles.Item(\"BuiltInTok\")\n$builtInTok.Font.Color = \"#008000\"\n"}
